$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$updates = @(
    @('D2', '311.69'),
    @('E2', '-4.57%'),
    @('D3', '41.07'),
    @('E3', '-6.94%'),
    @('D4', '5.108'),
    @('E4', '-3.48%'),
    @('D5', '0.07871'),
    @('E5', '-5.65%'),
    @('D6', '4.341'),
    @('E6', '-1.48%'),
    @('D7', '1.687'),
    @('E7', '-13.00%'),
    @('D8', '0.9225'),
    @('E8', '-4.76%'),
    @('D9', '0.1094'),
    @('E9', '-2.72%'),
    @('D10', '0.1780'),
    @('E10', '-5.81%'),
    @('D11', '0.09113'),
    @('E11', '-5.70%'),
    @('D12', '0.04396'),
    @('E12', '-4.37%'),
    @('D13', '7.188'),
    @('E13', '-15.86%'),
    @('D14', '0.1058'),
    @('E14', '-0.07%'),
    @('D15', '0.001261'),
    @('E15', '-2.37%'),
    @('D16', '0.005948'),
    @('E16', '1.30%'),
    @('B17', 'LEO'),
    @('C17', 'https://coinranking.com/coin/mqtUpyBxu8O8+leo-leo'),
    @('D17', '3.378'),
    @('E17', '-0.72%'),
    @('B18', 'BTSEToken'),
    @('C18', 'https://coinranking.com/coin/EOSL_JJKNMEr+btsetoken-btse'),
    @('D18', '2.559'),
    @('E18', '0.71%'),
    @('B19', 'BitpandaEcosystemToken'),
    @('C19', 'https://coinranking.com/coin/Uzf_Wjqc+bitpandaecosystemtoken-best'),
    @('D19', '0.3369'),
    @('E19', '0.35%'),
    @('B20', 'ProBitToken'),
    @('C20', 'https://coinranking.com/coin/lQP4d6T2+probittoken-prob'),
    @('D20', '0.1375'),
    @('E20', '0.25%'),
    @('B21', 'ZBToken'),
    @('C21', 'https://coinranking.com/coin/CxmvOsCyENPso+zbtoken-zb'),
    @('D21', '0.2800'),
    @('E21', '8.65%'),
    @('B22', 'CoinExToken'),
    @('C22', 'https://coinranking.com/coin/APDVU0XEViZ2o+coinextoken-cet'),
    @('D22', '0.04161'),
    @('E22', '0.07%'),
    @('B23', 'BitKan'),
    @('C23', 'https://coinranking.com/coin/RDOsLDgvY-AXe+bitkan-kan'),
    @('D23', '0.001222'),
    @('E23', '-0.93%'),
    @('B24', 'HotbitToken'),
    @('C24', 'https://coinranking.com/coin/uQJB8Ocu8lTb+hotbittoken-htb'),
    @('D24', '0.004153'),
    @('E24', '-5.81%'),
    @('D25', '0.0001225'),
    @('E25', '-5.91%'),
    @('D26', '0.0002990'),
    @('E26', '0.21%'),
    @('D38', '0.02457'),
    @('E38', '-8.22%'),
    @('D39', '0.05315'),
    @('E39', '-4.42%'),
    @('D40', '0.008008'),
    @('E40', '2.04%'),
    @('D41', '0.1355'),
    @('E41', '-3.93%'),
    @('D42', '0.007534'),
    @('E42', '2.85%'),
    @('D43', '0.001989'),
    @('E43', '-5.87%'),
    @('D44', '0.008192'),
    @('E44', '4.53%'),
    @('D45', '0.3105'),
    @('E45', '-11.38%'),
    @('D46', '0.00006766'),
    @('E46', '-1.36%'),
    @('D47', '0.00000000753'),
    @('E47', '0.21%'),
    @('D48', '0.003428'),
    @('E48', '-1.99%'),
    @('D49', '0.004116'),
    @('E49', '16.39%'),
    @('D50', '0.00002108'),
    @('E50', '0.21%'),
    @('D51', '0.0002008'),
    @('E51', '0.21%'),
)

foreach ($u in $updates) {
    $cellAddr = $u[0]
    $newVal = $u[1]
    $rng = $ws.Range($cellAddr)
    $origStyle = $rng.Style
    $rng.Value = "'" + $newVal
    $rng.Style = $origStyle
}
